# Applies the "Minor updates to data being read in, etc." revision to the
# Bull Trout Habitat Quality (RESTORATION) sheet:
#   - Entiat River Lake 04 / Potato 07 rows get revised Riparian scores & sums
#   - A new reach "Nason Creek Lower 03" is inserted as row 5, so every
#     subsequent "Nason Creek Lower NN" reach is renumbered down by one
#   - A new reach "Nason Creek Lower 14" is inserted as row 16, pushing the
#     former row 16 ("Nason Creek Lower 15") down to row 17
#   - Various individual attribute scores / derived sum+pct/flag columns are
#     corrected on several rows along the way

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Entiat River Lake 04 (HQ figures recomputed: Riparian-CanopyCover P2 & Riparian_Mean R2 revised)
$ws.Range("P2").Value = 3
$ws.Range("R2").Value = 4
$ws.Range("T2").Value = 32
$ws.Range("U2").Value = 0.7111111111111111

# Row 3: Entiat River Potato 07 (renamed attribute lists now include "Riparian"; scores + protection flag revised)
$ws.Range("P3").Value = 3
$ws.Range("R3").Value = 3
$ws.Range("T3").Value = 31
$ws.Range("U3").Value = 0.6888888888888889
$ws.Range("W3").Value = 1
$ws.Range("Y3").Value = "Stability,Cover-Wood,Flow-SummerBaseFlow,Off-Channel-Floodplain,Off-Channel-Side-Channels,Riparian,Temperature-Rearing"
$ws.Range("Z3").Value = "Stability,Cover-Wood,Flow-SummerBaseFlow,Off-Channel-Floodplain,Off-Channel-Side-Channels,Riparian,Temperature-Rearing"

# Row 5: Row 5 becomes "Nason Creek Lower 03" (new survey data for this reach)
$ws.Range("A5").Value = "Nason Creek Lower 03"
$ws.Range("H5").Value = 5
$ws.Range("I5").Value = 4
$ws.Range("K5").Value = 5
$ws.Range("N5").Value = 5
$ws.Range("O5").Value = 5
$ws.Range("T5").Value = 34
$ws.Range("U5").Value = 0.7555555555555555
$ws.Range("W5").Value = 3
$ws.Range("X5").Value = "Temperature-Rearing"
$ws.Range("Y5").Value = "Flow-SummerBaseFlow,Off-Channel-Floodplain,Riparian"
$ws.Range("Z5").Value = "Flow-SummerBaseFlow,Off-Channel-Floodplain,Riparian,Temperature-Rearing"

# Row 6: Row 6 relabeled "Nason Creek Lower 05" -> "Nason Creek Lower 04" (values unchanged)
$ws.Range("A6").Value = "Nason Creek Lower 04"

# Row 7: Row 7 relabeled -> "Nason Creek Lower 05"; Riparian scores revised
$ws.Range("A7").Value = "Nason Creek Lower 05"
$ws.Range("P7").Value = 1
$ws.Range("R7").Value = 2
$ws.Range("T7").Value = 24
$ws.Range("U7").Value = 0.5333333333333333

# Row 8: Row 8 relabeled -> "Nason Creek Lower 06" (values unchanged)
$ws.Range("A8").Value = "Nason Creek Lower 06"

# Row 9: Row 9 relabeled -> "Nason Creek Lower 07"; several scores revised, ChannelStability (H9) now populated
$ws.Range("A9").Value = "Nason Creek Lower 07"
$ws.Range("H9").Value = 3
$ws.Range("M9").Value = 3
$ws.Range("N9").Value = 3
$ws.Range("Q9").Value = 3
$ws.Range("R9").Value = 3
$ws.Range("T9").Value = 25
$ws.Range("U9").Value = 0.5555555555555556
$ws.Range("X9").Value = "PoolQuantity&Quality,Temperature-Rearing"
$ws.Range("Y9").Value = "Stability,Cover-Wood,Flow-SummerBaseFlow,Off-Channel-Floodplain,Off-Channel-Side-Channels,Riparian"

# Row 10: Row 10 relabeled -> "Nason Creek Lower 08"; BankStability/ChannelStability merge (H10 removed) + score revisions
$ws.Range("A10").Value = "Nason Creek Lower 08"
$ws.Range("G10").Value = 3
$ws.Range("H10").Value = $null
$ws.Range("I10").Value = 3
$ws.Range("K10").Value = 3
$ws.Range("O10").Value = 1
$ws.Range("P10").Value = 3
$ws.Range("Q10").Value = 1
$ws.Range("T10").Value = 20
$ws.Range("U10").Value = 0.4444444444444444
$ws.Range("X10").Value = "Off-Channel-Floodplain,Off-Channel-Side-Channels,PoolQuantity&Quality,Temperature-Rearing"
$ws.Range("Y10").Value = "Stability,Cover-Wood,Flow-SummerBaseFlow,Riparian"

# Row 11: Row 11 relabeled -> "Nason Creek Lower 09" (values unchanged)
$ws.Range("A11").Value = "Nason Creek Lower 09"

# Row 12: Row 12 relabeled -> "Nason Creek Lower 10"; many scores revised
$ws.Range("A12").Value = "Nason Creek Lower 10"
$ws.Range("G12").Value = 1
$ws.Range("I12").Value = 2
$ws.Range("J12").Value = 5
$ws.Range("K12").Value = 1
$ws.Range("M12").Value = 1
$ws.Range("O12").Value = 3
$ws.Range("P12").Value = 1
$ws.Range("R12").Value = 1
$ws.Range("T12").Value = 18
$ws.Range("U12").Value = 0.4
$ws.Range("X12").Value = "Cover-Wood,Off-Channel-Floodplain,Off-Channel-Side-Channels,Riparian,Temperature-Rearing"
$ws.Range("Y12").Value = "Stability,Flow-SummerBaseFlow,PoolQuantity&Quality"
$ws.Range("Z12").Value = "Stability,Cover-Wood,Flow-SummerBaseFlow,Off-Channel-Floodplain,Off-Channel-Side-Channels,PoolQuantity&Quality,Riparian,Temperature-Rearing"

# Row 13: Row 13 relabeled -> "Nason Creek Lower 11"; several scores revised
$ws.Range("A13").Value = "Nason Creek Lower 11"
$ws.Range("K13").Value = 3
$ws.Range("N13").Value = 1
$ws.Range("O13").Value = 1
$ws.Range("P13").Value = 3
$ws.Range("R13").Value = 2
$ws.Range("T13").Value = 20
$ws.Range("U13").Value = 0.4444444444444444
$ws.Range("X13").Value = "Off-Channel-Side-Channels,PoolQuantity&Quality,Temperature-Rearing"
$ws.Range("Y13").Value = "Stability,CoarseSubstrate,Cover-Wood,Flow-SummerBaseFlow,Off-Channel-Floodplain,Riparian"

# Row 14: Row 14 relabeled -> "Nason Creek Lower 12" (values unchanged)
$ws.Range("A14").Value = "Nason Creek Lower 12"

# Row 15: Row 15 relabeled -> "Nason Creek Lower 13" (values unchanged)
$ws.Range("A15").Value = "Nason Creek Lower 13"

# Row 16: Row 16 becomes "Nason Creek Lower 14" (new data, was formerly part of old row 16 "Lower 15")
$ws.Range("A16").Value = "Nason Creek Lower 14"
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 3
$ws.Range("P16").Value = 1
$ws.Range("Q16").Value = 1
$ws.Range("R16").Value = 1
$ws.Range("T16").Value = 21
$ws.Range("U16").Value = 0.4666666666666667
$ws.Range("X16").Value = "Cover-Wood,Riparian,Temperature-Rearing"
$ws.Range("Y16").Value = "Stability,CoarseSubstrate,Flow-SummerBaseFlow,Off-Channel-Floodplain,Off-Channel-Side-Channels,PoolQuantity&Quality"
$ws.Range("Z16").Value = "Stability,CoarseSubstrate,Cover-Wood,Flow-SummerBaseFlow,Off-Channel-Floodplain,Off-Channel-Side-Channels,PoolQuantity&Quality,Riparian,Temperature-Rearing"

# Row 17: New row 17 "Nason Creek Lower 15" (old row-16 data moved down one row)
$ws.Range("A17").Value = "Nason Creek Lower 15"
$ws.Range("B17").Value = "Wenatchee"
$ws.Range("C17").Value = "Lower Nason Creek"
$ws.Range("D17").Value = "yes"
$ws.Range("E17").Value = "yes"
$ws.Range("F17").Value = "yes"
$ws.Range("G17").Value = 3
$ws.Range("H17").Value = 3
$ws.Range("I17").Value = 3
$ws.Range("J17").Value = 3
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 5
$ws.Range("M17").Value = 3
$ws.Range("N17").Value = 3
$ws.Range("O17").Value = 3
$ws.Range("Q17").Value = 3
$ws.Range("R17").Value = 3
$ws.Range("S17").Value = 1
$ws.Range("T17").Value = 27
$ws.Range("U17").Value = 0.6
$ws.Range("V17").Value = 5
$ws.Range("W17").Value = 1
$ws.Range("X17").Value = "Temperature-Rearing"
$ws.Range("Y17").Value = "Stability,CoarseSubstrate,Cover-Wood,Off-Channel-Floodplain,Off-Channel-Side-Channels,PoolQuantity&Quality,Riparian"
$ws.Range("Z17").Value = "Stability,CoarseSubstrate,Cover-Wood,Off-Channel-Floodplain,Off-Channel-Side-Channels,PoolQuantity&Quality,Riparian,Temperature-Rearing"
